$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Diferencas")
$ws1.Range("B3").Value = -0.02821568627450709
$ws1.Range("C3").Value = 0.004477647058830195
$ws1.Range("D3").Value = -0.005736764705878084
$ws1.Range("B5").Value = -0.06277000000003408
$ws1.Range("C5").Value = -0.05416999999991912
$ws1.Range("D5").Value = -0.05342999999996845
$ws1.Range("B6").Value = -0.03916249999995169
$ws1.Range("C6").Value = -0.007815000000027328
$ws1.Range("D6").Value = -0.01580999999996446
$ws1.Range("B7").Value = -0.02993000000001256
$ws1.Range("C7").Value = 0.007289999999978147
$ws1.Range("D7").Value = -0.004869999999959407
$ws1.Range("B8").Value = -0.0172275000000236
$ws1.Range("C8").Value = 0.02247000000005006
$ws1.Range("D8").Value = 0.003135000000009325
$ws1.Range("B9").Value = 0.01482999999992718
$ws1.Range("C9").Value = 0.04570999999995184
$ws1.Range("D9").Value = 0.05002999999995883

$ws2 = $wb.Worksheets.Item("Diferencas_Absolutas")
$ws2.Range("B3").Value = 0.02872666666666189
$ws2.Range("C3").Value = 0.01974058823528968
$ws2.Range("D3").Value = 0.0134865686274393
$ws2.Range("B4").Value = 0.01432565781900236
$ws2.Range("C4").Value = 0.01324251477893625
$ws2.Range("D4").Value = 0.01171476875234177
$ws2.Range("B5").Value = 0.0008099999999212626
$ws2.Range("C5").Value = 0.0003000000000047409
$ws2.Range("D5").Value = 0.0001299999999792023
$ws2.Range("B6").Value = 0.0172275000000236
$ws2.Range("C6").Value = 0.007604999999989093
$ws2.Range("D6").Value = 0.004259999999936315
$ws2.Range("B7").Value = 0.02993000000001256
$ws2.Range("C7").Value = 0.01735499999995227
$ws2.Range("D7").Value = 0.01014500000002272
$ws2.Range("B8").Value = 0.03916249999995169
$ws2.Range("C8").Value = 0.02915999999998509
$ws2.Range("D8").Value = 0.0188424999999936
$ws2.Range("B9").Value = 0.06277000000003408
$ws2.Range("C9").Value = 0.05416999999991912
$ws2.Range("D9").Value = 0.05342999999996845
